$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old layout had a label column (A) plus six weighting columns
# (Ödev1, Ödev2, Quiz, Quiz4, Vize, Final) and a Final "Toplam" sum column (H).
# The new layout keeps only the label column (A) and a single "Toplam" column (B)
# that sums each row. Deleting the entire columns C:H shifts H (and its formula)
# left into B and drops the now-unused weighting values, renumbering the shared
# strings automatically.
$ws.Range("C1:H1").EntireColumn.Delete()

# The surviving column B used to be headed "Ödev1"; rename it to "Toplam".
$ws.Range("B1").Value = "Toplam"

# Recreate each row's total as a formula over the row's two remaining cells.
for ($row = 2; $row -le 6; $row++) {
    $cellRef = "B" + $row
    $ws.Range($cellRef).Formula = "=SUM(B" + $row + ":A" + $row + ")"
}
